$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 9: simplify the merged "Raca" header (E9:G9) border so it matches
# the plain bordered/centered header style already used by "TipoPet" (I9:J9). ---
$ws.Range("I9:J9").Copy()
$ws.Range("E9:G9").PasteSpecial(-4122)

# --- Fix breed-name typo: "Shi-tzu" -> "Shih-tzu" ---
$ws.Range("F11").Value2 = "Shih-tzu"

# --- Format the new data rows/cells to match the existing bordered table look
# before dropping in their values (order matters for shared-string append order). ---

# Text cells (nomeRaca / NomeTipoPet) - match F13's plain bordered text style.
$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("J13").PasteSpecial(-4122)

# Numeric id cells - match E13's bordered/centered numeric style.
$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("I13").PasteSpecial(-4122)

# --- Values: new breed rows (Raca table, E:G) ---
$ws.Range("F14").Value2 = "Sphynx"
$ws.Range("F15").Value2 = "Persa"
$ws.Range("F16").Value2 = "Mangalarga"

$ws.Range("E14").Value2 = 4
$ws.Range("G14").Value2 = 2
$ws.Range("E15").Value2 = 5
$ws.Range("G15").Value2 = 2
$ws.Range("E16").Value2 = 6
$ws.Range("G16").Value2 = 3

# --- Values: fix idTipoPet for Vira-Lata (dog, not cat) + new TipoPet row (Cavalo) ---
$ws.Range("G13").Value2 = 1
$ws.Range("I13").Value2 = 3
$ws.Range("J13").Value2 = "Cavalo"

# --- Sheet view: scrolled right with a new active selection ---
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("O9").Select()
